# The "Requisitos" bullet list ends with three lines separated by manual
# line breaks (w:br). Originally the order is:
#   LOT2053 - Microbiologia (Requisito fraco)
#   LOT2007 - Bioquímica I  (Requisito fraco)
#   LOT2040 - Engenharia Genética (Requisito fraco)
# and it needs to become:
#   LOT2007 - Bioquímica I  (Requisito fraco)
#   LOT2040 - Engenharia Genética (Requisito fraco)
#   LOT2053 - Microbiologia (Requisito fraco)
# i.e. the "Microbiologia" line moves from first to last place.

$d = $word.ActiveDocument
$lb = [char]11   # manual line break character produced by <w:br/>

$micro = "LOT2053 -  Microbiologia  (Requisito fraco)"
$bioq  = "LOT2007 -  Bioquímica I  (Requisito fraco)"
$eng   = "LOT2040 -  Engenharia Genética  (Requisito fraco)"

# Find the start of the "Microbiologia" line (the first line of the list).
$find = $d.Content.Find
$find.ClearFormatting()
$found = $find.Execute($micro + $lb, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find the LOT2053 Microbiologia requirement line."
}
$anchor = $find.Parent.Start

# Insert fresh copies of the Bioquímica/Engenharia lines immediately before
# that anchor (in reverse order, since each insert pushes the later one
# further along) so they become new, distinct runs placed ahead of the
# original Microbiologia run.
$insPoint = $d.Range($anchor, $anchor)
$insPoint.InsertBefore($eng + $lb)

$insPoint2 = $d.Range($anchor, $anchor)
$insPoint2.InsertBefore($bioq + $lb)

# The paragraph now reads: Bioquímica, Engenharia, Microbiologia, Bioquímica,
# Engenharia (original two trailing lines are now duplicated at the end).
# Locate that whole sequence and drop the now-redundant trailing pair,
# leaving only the freshly reordered three lines.
$find2 = $d.Content.Find
$find2.ClearFormatting()
$found2 = $find2.Execute($micro + $lb + $bioq + $lb + $eng + $lb, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found2) {
    throw "Could not find the duplicated trailing requirement lines."
}
$tailStart = $find2.Parent.Start + ($micro + $lb).Length
$tailEnd = $find2.Parent.End
$tail = $d.Range($tailStart, $tailEnd)
$tail.Delete()
